$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 12 (Training -> Save button) ---
$ws.Range("C12").Value = "assert,click,click"
$ws.Range("D12").Value = "Trainings,no value,no value"
$ws.Range("E12").Value = "training_assert,click_training,finish"
$ws.Range("F12").Value = "yes"

# --- New row 13 ---
$ws.Range("A13").Value = "141289"
$ws.Range("B13").Value = "Client->Client Summary->Training->To verify that the user can view previously entered training data and add new entries on the Client Record-Training Data screen."
$ws.Range("C13").Value = "click,date,type,click"
$ws.Range("D13").Value = "no value,faker,faker,no value"
$ws.Range("E13").Value = "next,date,first_name_next,cancel"
$ws.Range("F13").Value = "yes"

# --- New row 14 ---
$ws.Range("A14").Value = "141370"
$ws.Range("B14").Value = "To verify that the user can upload instructions by clicking the upload button."
$ws.Range("C14").Value = "click,click,click"
$ws.Range("D14").Value = "no value,no value,no value"
$ws.Range("E14").Value = "instruction,next,close_popup"
$ws.Range("F14").Value = "yes"

# --- New row 15 ---
$ws.Range("A15").Value = "144229"
$ws.Range("B15").Value = "To verify that the existing summary page displays a detailed view of all the tabs in the Summary Page."
$ws.Range("C15").Value = "click"
$ws.Range("D15").Value = "no value"
$ws.Range("E15").Value = "training"
$ws.Range("F15").Value = "no"

# --- New row 16 ---
$ws.Range("A16").Value = "141406"
$ws.Range("B16").Value = "Client-Client Summary-Activity log->To verify if the correct action of the user mentioned in the FRD (Record updated, Case manager assigned, Status changed, Referral created, comments added, documents uploaded) is displayed in the Activity log."
$ws.Range("C16").Value = "click"
$ws.Range("D16").Value = "no value"
$ws.Range("E16").Value = "document"
$ws.Range("F16").Value = "yes"

# --- New row 17 ---
$ws.Range("A17").Value = "141340"
$ws.Range("B17").Value = "To verify that each document should have at least one tag in the Upload document pop-up."
$ws.Range("C17").Value = "click"
$ws.Range("D17").Value = "no value"
$ws.Range("E17").Value = "instruction"
$ws.Range("F17").Value = "yes"

# --- New row 18 ---
$ws.Range("A18").Value = "141370"
$ws.Range("B18").Value = "Client-Client Summary->Instruction->To verify that the user can upload instructions by clicking the upload button."
$ws.Range("C18").Value = "click"
$ws.Range("D18").Value = "no value"
$ws.Range("E18").Value = "training"
$ws.Range("F18").Value = "yes"

# --- New row 19 ---
$ws.Range("A19").Value = "141380"
$ws.Range("B19").Value = "Client->Client Summary->Instruction->To verify that if there is no data in the instruction section, the error message `"No record found`" is displayed."
$ws.Range("C19").Value = "click"
$ws.Range("D19").Value = "no value"
$ws.Range("E19").Value = "document"
$ws.Range("F19").Value = "yes"

# --- New row 20 ---
$ws.Range("A20").Value = "141406"
$ws.Range("B20").Value = "Client-Client Summary-Activity log->To verify if the correct action of the user mentioned in the FRD (Record updated, Case manager assigned, Status changed, Referral created, comments added, documents uploaded) is displayed in the Activity log."
$ws.Range("C20").Value = "click"
$ws.Range("D20").Value = "no value"
$ws.Range("E20").Value = "instruction"
$ws.Range("F20").Value = "yes"
